# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45190 (2023-09-21) to 45192 (2023-09-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 262 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45192
